# Weekly update: a new price observation is inserted at the top of the
# "Femacal de La Calera - Cilantro" data block (row 547), pushing the
# existing historical rows (547-650) down by one (548-651) and growing
# the used range from A1:R650 to A1:R651.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 547; Excel shifts 547:650 down to 548:651 and
# extends the sheet dimension automatically.
$ws.Rows("547:547").Insert()

# Populate the newly inserted row with the latest weekly record.
$ws.Range("A547").Value = 3
$ws.Range("B547").Value = "Femacal de La Calera"
$ws.Range("C547").Value = "Coquimbo"
$ws.Range("D547").Value = 45209
$ws.Range("E547").Value = 5
$ws.Range("F547").Value = 100112040
$ws.Range("G547").Value = "Cilantro"
$ws.Range("H547").Value = "Sin especificar"
$ws.Range("I547").Value = "Primera"
$ws.Range("J547").Value = 170
$ws.Range("K547").Value = 3500
$ws.Range("L547").Value = 4000
$ws.Range("M547").Value = 3765
$ws.Range("N547").Value = "`$/docena de atados (3 kilos)"
$ws.Range("O547").Value = "Provincia de Quillota"
$ws.Range("P547").Value = 1255
$ws.Range("Q547").Value = 3
$ws.Range("R547").Value = "Hortaliza"
